$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that lives on H52 (mailto:kingmahimn@gmail.com)
$ws.Hyperlinks.Delete()

# Clear out the last data row (row 52), which contained the workbook
# author's own personal sample record. The H52 cell keeps its
# (hyperlink) style but becomes empty, while the rest of the row's
# cells are removed entirely.
$ws.Range("A52:G52").ClearContents()
$ws.Range("H52").ClearContents()
$ws.Range("I52:J52").ClearContents()

# Update the active selection to match the author's last position
$ws.Range("M25").Select()
